# Trade #135 (MarketMaking) closes early with a small gain, and a brand new
# trade #168 is opened. This ripples through the Summary, Strategy Status,
# All Trades and MarketMaking sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1401.49   # Current Capital
$summary.Range("B4").Value = 1.28      # Total P&L $
$summary.Range("B6").Value = 135       # Total Trades
$summary.Range("B7").Value = 58        # Winning Trades
$summary.Range("B9").Value = 42.96     # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet - row 5 is the MarketMaking strategy
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 101.49     # Capital
$status.Range("D5").Value = 102        # Trades
$status.Range("E5").Value = 1.17       # P&L $
$status.Range("F5").Value = 1.49       # P&L %
$status.Range("G5").Value = 43.14      # Win Rate %

# ---------------------------------------------------------------------------
# All Trades sheet
#   columns: A Trade#, B Date, C Time, D Strategy, E Side, F Entry Price,
#            G Exit Price, H Status, I P&L %, J P&L $, K Capital After,
#            L Exit Reason, M Duration (min), N Entry Slippage,
#            O Exit Slippage, P Confidence, Q Entry Reason
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Row 136 = trade #135 : closes out
$allTrades.Cells.Item(136, 7).Value = 0.92             # G - Exit Price
$allTrades.Cells.Item(136, 8).Value = "CLOSED"         # H - Status
$allTrades.Cells.Item(136, 9).Value = 3.3708           # I - P&L %
$allTrades.Cells.Item(136, 10).Value = 0.03            # J - P&L $
$allTrades.Cells.Item(136, 11).Value = 101.49          # K - Capital After
$allTrades.Cells.Item(136, 12).Value = "early_exit"    # L - Exit Reason
$allTrades.Cells.Item(136, 13).Value = 0.15            # M - Duration (min)

# Row 169 = brand new trade #168, still OPEN. Date/time columns must stay
# text, so force the number format to Text before assigning so COM doesn't
# silently coerce the date-looking string into a date serial.
$allTrades.Cells.Item(169, 1).Value = 168                          # A - Trade #
$allTrades.Cells.Item(169, 2).NumberFormat = "@"
$allTrades.Cells.Item(169, 2).Value = "2026-02-17"                 # B - Date
$allTrades.Cells.Item(169, 3).NumberFormat = "@"
$allTrades.Cells.Item(169, 3).Value = "21:32:36"                   # C - Time
$allTrades.Cells.Item(169, 4).Value = "MarketMaking"                # D - Strategy
$allTrades.Cells.Item(169, 5).Value = "DOWN"                        # E - Side
$allTrades.Cells.Item(169, 6).Value = 0.89                          # F - Entry Price
$allTrades.Cells.Item(169, 8).Value = "OPEN"                        # H - Status
$allTrades.Cells.Item(169, 9).Value = 0                             # I - P&L %
$allTrades.Cells.Item(169, 10).Value = 0                            # J - P&L $
$allTrades.Cells.Item(169, 11).Value = 101.4641758035408            # K - Capital After
$allTrades.Cells.Item(169, 13).Value = 0                            # M - Duration (min)
$allTrades.Cells.Item(169, 14).Value = 0                            # N - Entry Slippage
$allTrades.Cells.Item(169, 15).Value = 0                            # O - Exit Slippage
$allTrades.Cells.Item(169, 16).Value = 0.6                          # P - Confidence
$allTrades.Cells.Item(169, 17).Value = "Normal spread capture: 19600 bps"  # Q - Entry Reason

# ---------------------------------------------------------------------------
# MarketMaking sheet
#   columns: A Trade#, B Date, C Time, D Strategy, E Side, F Entry Price,
#            G Exit Price, H Status, I P&L %, J P&L $, K Capital After,
#            L Entry Slippage, M Exit Slippage, N Confidence, O Entry Reason,
#            P Exit Reason, Q Duration (min)
# ---------------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

# Row 103 = trade #135 : closes out
$mm.Cells.Item(103, 7).Value = 0.92             # G - Exit Price
$mm.Cells.Item(103, 8).Value = "CLOSED"         # H - Status
$mm.Cells.Item(103, 9).Value = 3.3708           # I - P&L %
$mm.Cells.Item(103, 10).Value = 0.03            # J - P&L $
$mm.Cells.Item(103, 11).Value = 101.49          # K - Capital After
$mm.Cells.Item(103, 16).Value = "early_exit"    # P - Exit Reason
$mm.Cells.Item(103, 17).Value = 0.15            # Q - Duration (min)

# Row 136 = brand new trade #168, still OPEN.
$mm.Cells.Item(136, 1).Value = 168                          # A - Trade #
$mm.Cells.Item(136, 2).NumberFormat = "@"
$mm.Cells.Item(136, 2).Value = "2026-02-17"                 # B - Date
$mm.Cells.Item(136, 3).NumberFormat = "@"
$mm.Cells.Item(136, 3).Value = "21:32:36"                   # C - Time
$mm.Cells.Item(136, 4).Value = "MarketMaking"                # D - Strategy
$mm.Cells.Item(136, 5).Value = "DOWN"                        # E - Side
$mm.Cells.Item(136, 6).Value = 0.89                          # F - Entry Price
$mm.Cells.Item(136, 8).Value = "OPEN"                        # H - Status
$mm.Cells.Item(136, 9).Value = 0                             # I - P&L %
$mm.Cells.Item(136, 10).Value = 0                            # J - P&L $
$mm.Cells.Item(136, 11).Value = 101.4641758035408            # K - Capital After
$mm.Cells.Item(136, 12).Value = 0                            # L - Entry Slippage
$mm.Cells.Item(136, 13).Value = 0                            # M - Exit Slippage
$mm.Cells.Item(136, 14).Value = 0.6                          # N - Confidence
$mm.Cells.Item(136, 15).Value = "Normal spread capture: 19600 bps"  # O - Entry Reason
$mm.Cells.Item(136, 17).Value = 0                            # Q - Duration (min)

Write-Output "Applied trade #135 close + trade #168 open across all sheets"
